$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (shifts existing rows 17-36 down to 18-37)
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with this week's data
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44447
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112009
$ws.Range("G17").Value = "Acelga"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 1100
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = 1150
$ws.Range("N17").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 383
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = "Hortaliza"
